$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / title text updates: October -> November ---
$ws.Range("A1").Value = "Table 4.11.B. Average Cost of Petroleum Liquids Delivered for Electricity Generation by State, (Year-to-Date) November 2016 and 2015"
$ws.Range("B4").Value = "November 2016 YTD"
$ws.Range("C4").Value = "November 2015 YTD"
$ws.Range("E4").Value = "November 2016 YTD"
$ws.Range("F4").Value = "November 2015 YTD"
$ws.Range("G4").Value = "November 2016 YTD"
$ws.Range("H4").Value = "November 2015 YTD"

# --- Data cell updates (row by row) ---
# Row 5
$ws.Range("B5").Value = "W"
$ws.Range("E5").Value = 9.64
$ws.Range("F5").Value = 11.49
$ws.Range("G5").Value = "W"
# Row 6
$ws.Range("B6").Value = 10.59
$ws.Range("G6").Value = 10.59
# Row 8
$ws.Range("C8").Value = 11.85
$ws.Range("E8").Value = 9.61
$ws.Range("H8").Value = 11.81
# Row 9
$ws.Range("B9").Value = 9.85
$ws.Range("E9").Value = 9.85
$ws.Range("F9").Value = 10.54
# Row 12
$ws.Range("B12").Value = 10.02
$ws.Range("C12").Value = 11.44
$ws.Range("D12").Value = -0.12
$ws.Range("E12").Value = 7.8
$ws.Range("F12").Value = 8.98
$ws.Range("H12").Value = 13.12
# Row 13
$ws.Range("B13").Value = 9.83
$ws.Range("C13").Value = 14.01
$ws.Range("D13").Value = -0.3
$ws.Range("G13").Value = 9.83
$ws.Range("H13").Value = 14.01
# Row 14
$ws.Range("B14").Value = 9.81
$ws.Range("C14").Value = 10.93
$ws.Range("E14").Value = 7.8
$ws.Range("F14").Value = 8.98
$ws.Range("G14").Value = 11.89
$ws.Range("H14").Value = 12.94
# Row 15
$ws.Range("B15").Value = 10.39
$ws.Range("C15").Value = 13.37
$ws.Range("D15").Value = -0.22
$ws.Range("G15").Value = 10.39
$ws.Range("H15").Value = 13.37
# Row 16
$ws.Range("C16").Value = 13.87
$ws.Range("E16").Value = 10.64
$ws.Range("F16").Value = 13.91
$ws.Range("H16").Value = 13.78
# Row 17
$ws.Range("B17").Value = 10.74
$ws.Range("F17").Value = 14.2
$ws.Range("G17").Value = 10.74
# Row 18
$ws.Range("B18").Value = 10.43
$ws.Range("C18").Value = 14.02
$ws.Range("D18").Value = -0.26
$ws.Range("E18").Value = 10.43
$ws.Range("F18").Value = 14.02
# Row 19
$ws.Range("B19").Value = 10.44
$ws.Range("C19").Value = 13.18
$ws.Range("D19").Value = -0.21
$ws.Range("E19").Value = 10.44
$ws.Range("F19").Value = 13.18
# Row 20
$ws.Range("C20").Value = 13.8
$ws.Range("E20").Value = 10.9
$ws.Range("F20").Value = 13.95
$ws.Range("H20").Value = 13.7
# Row 21
$ws.Range("E21").Value = 11.39
$ws.Range("F21").Value = 15.1
# Row 22
$ws.Range("C22").Value = 12.82
$ws.Range("E22").Value = 10.49
$ws.Range("F22").Value = 12.82
# Row 23
$ws.Range("B23").Value = 11.08
$ws.Range("C23").Value = 13.02
$ws.Range("E23").Value = 11.08
$ws.Range("F23").Value = 13.02
# Row 24
$ws.Range("B24").Value = 10.3
$ws.Range("C24").Value = 12.67
$ws.Range("D24").Value = -0.19
$ws.Range("E24").Value = 10.3
$ws.Range("F24").Value = 12.67
# Row 25
$ws.Range("C25").Value = 13.52
$ws.Range("E25").Value = 11.29
$ws.Range("F25").Value = 13.52
# Row 26
$ws.Range("B26").Value = 10.57
$ws.Range("C26").Value = 13.19
$ws.Range("D26").Value = -0.2
$ws.Range("E26").Value = 10.57
$ws.Range("F26").Value = 13.19
# Row 27
$ws.Range("C27").Value = 20.52
$ws.Range("D27").Value = -0.46
$ws.Range("F27").Value = 20.52
# Row 28
$ws.Range("B28").Value = 9.12
$ws.Range("C28").Value = 12.83
$ws.Range("D28").Value = -0.29
$ws.Range("E28").Value = 9.12
$ws.Range("F28").Value = 12.83
# Row 29
$ws.Range("B29").Value = 8.54
$ws.Range("D29").Value = -0.099
$ws.Range("E29").Value = 8.54
# Row 30
$ws.Range("B30").Value = 9.89
$ws.Range("C30").Value = 12.89
$ws.Range("E30").Value = 9.78
$ws.Range("F30").Value = 12.65
$ws.Range("G30").Value = 10.39
$ws.Range("H30").Value = 13.7
# Row 33
$ws.Range("F33").Value = 14.5
# Row 34
$ws.Range("B34").Value = 9.25
$ws.Range("C34").Value = 16.27
$ws.Range("D34").Value = -0.43
$ws.Range("E34").Value = 9.54
$ws.Range("F34").Value = 17.28
$ws.Range("G34").Value = 7.83
# Row 35
$ws.Range("B35").Value = 9.62
$ws.Range("D35").Value = -0.11
$ws.Range("G35").Value = 9.62
# Row 36
$ws.Range("E36").Value = 10.1
$ws.Range("F36").Value = 13.35
# Row 37
$ws.Range("B37").Value = 10.99
$ws.Range("C37").Value = 15.02
$ws.Range("D37").Value = -0.27
$ws.Range("E37").Value = 10.99
$ws.Range("F37").Value = 15.02
# Row 38
$ws.Range("C38").Value = 11.93
$ws.Range("E38").Value = 7.97
$ws.Range("H38").Value = 17.62
# Row 39
$ws.Range("B39").Value = "W"
$ws.Range("E39").Value = 11.22
$ws.Range("F39").Value = 14.03
$ws.Range("G39").Value = "W"
# Row 40
$ws.Range("E40").Value = 10.27
$ws.Range("F40").Value = 12.92
# Row 41
$ws.Range("E41").Value = 9.81
$ws.Range("F41").Value = 13.29
# Row 42
$ws.Range("B42").Value = 10.47
$ws.Range("C42").Value = 13.67
$ws.Range("D42").Value = -0.23
$ws.Range("E42").Value = 10.47
$ws.Range("F42").Value = 13.67
# Row 43
$ws.Range("B43").Value = 9.31
$ws.Range("D43").Value = -0.11
$ws.Range("E43").Value = 9.31
# Row 44
$ws.Range("B44").Value = 10.32
$ws.Range("C44").Value = 12.66
$ws.Range("D44").Value = -0.18
$ws.Range("E44").Value = 10.32
$ws.Range("F44").Value = 12.66
# Row 45
$ws.Range("B45").Value = 10.51
$ws.Range("C45").Value = 13.21
$ws.Range("D45").Value = -0.2
$ws.Range("E45").Value = 10.33
$ws.Range("F45").Value = 13.09
$ws.Range("G45").Value = 11.04
$ws.Range("H45").Value = 13.49
# Row 46
$ws.Range("E46").Value = 9.99
$ws.Range("F46").Value = 13.26
# Row 49
$ws.Range("E49").Value = 10.31
$ws.Range("F49").Value = 13.37
# Row 50
$ws.Range("C50").Value = 14.89
$ws.Range("E50").Value = 11.21
$ws.Range("F50").Value = 14.92
$ws.Range("H50").Value = 14.56
# Row 51
$ws.Range("B51").Value = 11.19
$ws.Range("C51").Value = 13.95
$ws.Range("D51").Value = -0.2
$ws.Range("E51").Value = 11.19
$ws.Range("F51").Value = 13.95
# Row 55
$ws.Range("F55").Value = 16.98
# Row 56
$ws.Range("B56").Value = 10.96
$ws.Range("C56").Value = 15.74
$ws.Range("D56").Value = -0.3
$ws.Range("E56").Value = 10.96
$ws.Range("F56").Value = 15.74
# Row 57
$ws.Range("B57").Value = 11.64
$ws.Range("E57").Value = 11.64
$ws.Range("F57").Value = 14.84
# Row 58
$ws.Range("B58").Value = 11.4
$ws.Range("C58").Value = 14.44
$ws.Range("D58").Value = -0.21
$ws.Range("E58").Value = 11.4
$ws.Range("F58").Value = 14.44
# Row 59
$ws.Range("E59").Value = 12.29
$ws.Range("F59").Value = 11.66
# Row 61
$ws.Range("C61").Value = 11.72
$ws.Range("F61").Value = 11.72
# Row 62
$ws.Range("E62").Value = 12.29
$ws.Range("F62").Value = 11.29
# Row 63
$ws.Range("E63").Value = 8.39
$ws.Range("F63").Value = 11.11
# Row 64
$ws.Range("B64").Value = 14.31
$ws.Range("C64").Value = 17.45
$ws.Range("D64").Value = -0.18
$ws.Range("E64").Value = 14.31
$ws.Range("F64").Value = 17.45
# Row 65
$ws.Range("E65").Value = 8.38
$ws.Range("F65").Value = 11.1
# Row 66
$ws.Range("B66").Value = 9.24
$ws.Range("C66").Value = 11.65
$ws.Range("D66").Value = -0.21
$ws.Range("E66").Value = 9.03
$ws.Range("F66").Value = 11.57
$ws.Range("G66").Value = 9.85
$ws.Range("H66").Value = 11.78
